$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new journal entries (rows 26 and 27) ---
# Copy the formatting of the last filled-in row (25) down into rows 26 and 27
# so the cells get the same date/time/wrap-text styles used throughout the table.
$ws.Range("A25:E25").Copy()
$ws.Range("A26:E26").PasteSpecial(-4122)
$ws.Range("A27:E27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 26: "Suite de la réalisation de la documentation"
$ws.Range("A26").Value2 = 43930
$ws.Range("B26").Value2 = 0.4236111111111111
$ws.Range("C26").Value2 = 0.44097222222222227
$ws.Range("E26").Value2 = "Suite de la réalisation de la documentation"

# Row 27: "Revue finale de la documentation avant le rendu"
$ws.Range("A27").Value2 = 43930
$ws.Range("B27").Value2 = 0.63888888888888895
$ws.Range("C27").Value2 = 0.70486111111111116
$ws.Range("E27").Value2 = "Revue finale de la documentation avant le rendu"

# Extend the "duration" formula (C-B) from D2:D25 down through D27
$ws.Range("D2:D27").Formula = "=C2-B2"

# --- Update the active selection to D31 (matches the saved cursor position) ---
$ws.Range("D31").Select()
